# Daily attendance processing - 2025-11-16 05:46:42
# Swap the order of names in the "Recorded By" column (G) for rows where
# exactly two comma-separated entries are recorded and neither of them is
# the backup account (backup@backdoor.com). Rows with a single entry, three
# entries, or involving the backup account are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text.Contains("backup@backdoor.com")) {
        continue
    }

    $parts = $text -split ", "
    if ($parts.Count -eq 2) {
        $newValue = $parts[1] + ", " + $parts[0]
        $cell.Value = $newValue
    }
}
